$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the default cell style from "Standard" (de-DE) to "Normal" (en-US) ---
# The engine only picks up a rename when the old builtin style is removed and
# the new name is (re-)added against the same builtin slot.
$wb.Styles.Item(1).Delete()
$wb.Styles.Add("Normal")

# --- Extend the time log with the two new rows for the final day ---
# Row 59
$ws.Range("A59").Value = "14.4.2020"
$ws.Range("A59").HorizontalAlignment = -4152
$ws.Range("A59").VerticalAlignment = -4108

$ws.Range("B59").Value = 0.375
$ws.Range("B59").NumberFormat = "h:mm"

$ws.Range("C59").Value = 0.40625
$ws.Range("C59").NumberFormat = "h:mm"

$ws.Range("E59").Value = "Project"
$ws.Range("E59").HorizontalAlignment = -4131
$ws.Range("E59").VerticalAlignment = -4108

$ws.Range("F59").Value = "Documentation"
$ws.Range("F59").HorizontalAlignment = -4131
$ws.Range("F59").VerticalAlignment = -4108

$ws.Range("G59").Value = "Last documentation of files"

# Row 60
$ws.Range("A60").Value = "14.4.2020"
$ws.Range("A60").HorizontalAlignment = -4152
$ws.Range("A60").VerticalAlignment = -4108

$ws.Range("B60").Value = 0.03125
$ws.Range("B60").NumberFormat = "h:mm"

$ws.Range("C60").Value = 0.42708333333333331
$ws.Range("C60").NumberFormat = "h:mm"

$ws.Range("E60").Value = "Project"
$ws.Range("E60").HorizontalAlignment = -4131
$ws.Range("E60").VerticalAlignment = -4108

$ws.Range("F60").Value = "Documentation"
$ws.Range("F60").HorizontalAlignment = -4131
$ws.Range("F60").VerticalAlignment = -4108

$ws.Range("G60").Value = "Write documentation PDF"

# Fill the D column (elapsed time) for the two new rows as one shared formula,
# matching the existing D5:D58 shared-formula pattern used throughout the sheet.
$ws.Range("D59:D60").Formula = "=C59-B59"
$ws.Range("D59:D60").NumberFormat = "h:mm"
$ws.Range("D59:D60").HorizontalAlignment = -4152
$ws.Range("D59:D60").VerticalAlignment = -4108

# Move the active selection to reflect the new last row of the log.
$ws.Range("G61").Select()
